# Applies FPVA sheet updates: exit pressure (F) and supersonic volume % (G) samples,
# plus "more..." annotations on two rows (no min/max yet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FPVA")

$ws.Range("F3").Value = 2.379025
$ws.Range("G3").Value = 98.87363

$ws.Range("F4").Value = 2.389668
$ws.Range("G4").Value = 98.84409

$ws.Range("F5").Value = 2.393329
$ws.Range("G5").Value = 98.82307

$ws.Range("F6").Value = 2.402416
$ws.Range("G6").Value = 98.80643

$ws.Range("F7").Value = 2.399836
$ws.Range("G7").Value = 98.80103

$ws.Range("F8").Value = 2.386293
$ws.Range("G8").Value = 98.84968

$ws.Range("F9").Value = 2.388704
$ws.Range("G9").Value = 98.83205

$ws.Range("F10").Value = 2.389445
$ws.Range("G10").Value = 98.84048

$ws.Range("F11").Value = 2.385821
$ws.Range("G11").Value = 98.85662

$ws.Range("F12").Value = 2.403574
$ws.Range("G12").Value = 98.77254

$ws.Range("F13").Value = 2.393145
$ws.Range("G13").Value = 98.81339

$ws.Range("F14").Value = 2.380624
$ws.Range("G14").Value = 98.86236

$ws.Range("F15").Value = 2.400453
$ws.Range("G15").Value = 98.79408

$ws.Range("F16").Value = 2.394499
$ws.Range("G16").Value = 98.82126
$ws.Range("H16").Value = "more…"

$ws.Range("F17").Value = 2.388582
$ws.Range("G17").Value = 98.82803

$ws.Range("F18").Value = 2.408596
$ws.Range("G18").Value = 98.77375

$ws.Range("F19").Value = 2.405758
$ws.Range("G19").Value = 98.7832

$ws.Range("F20").Value = 2.401817
$ws.Range("G20").Value = 98.78914

$ws.Range("F21").Value = 2.398059
$ws.Range("G21").Value = 98.79341

$ws.Range("F22").Value = 2.397976
$ws.Range("G22").Value = 98.81894

$ws.Range("F23").Value = 2.400359
$ws.Range("G23").Value = 98.80613

$ws.Range("F24").Value = 2.379926
$ws.Range("G24").Value = 98.86254

$ws.Range("F25").Value = 2.393892
$ws.Range("G25").Value = 98.82293

$ws.Range("F26").Value = 2.410452
$ws.Range("G26").Value = 98.75828

$ws.Range("F27").Value = 2.399798
$ws.Range("G27").Value = 98.80267
$ws.Range("H27").Value = "more…"

# Apply the 5-decimal number format to the whole supersonic-volume column (G3:G31),
# including the trailing rows that get a formatted-but-empty placeholder cell.
$ws.Range("G3:G31").NumberFormat = "0.00000"

# Update the sheet selection to match the authored state.
$ws.Activate()
$ws.Range("G27").Select()
